$wb = $excel.ActiveWorkbook

# --- Traits: delete data rows 2:5 (selection was B6, now A2:B5) ---
$ws = $wb.Worksheets.Item("Traits")
$ws.Activate()
$ws.Range("A2:B5").Select()
$excel.Selection.EntireRow.Delete()

# --- Soils: keep row 2's C2 style placeholder, drop A2 value; delete rows 3:5 ---
$ws = $wb.Worksheets.Item("Soils")
$ws.Activate()
$ws.Range("A3:A5").EntireRow.Delete()
$ws.Range("A2").Clear()

# --- SoilLayers: delete data rows 2:5 (selection was B5, now A2:B5) ---
$ws = $wb.Worksheets.Item("SoilLayers")
$ws.Activate()
$ws.Range("A2:B5").Select()
$excel.Selection.EntireRow.Delete()

# --- Regions: delete data rows 2:5 (selection was A5, now A2:A5) ---
$ws = $wb.Worksheets.Item("Regions")
$ws.Activate()
$ws.Range("A2:A5").Select()
$excel.Selection.EntireRow.Delete()

# --- Sites: delete data rows 2:9 (selection was B10, now A2:B9) ---
$ws = $wb.Worksheets.Item("Sites")
$ws.Activate()
$ws.Range("A2:B9").Select()
$excel.Selection.EntireRow.Delete()

# --- Fields: delete data rows 2:9 (selection was C10, now A2:C9) ---
$ws = $wb.Worksheets.Item("Fields")
$ws.Activate()
$ws.Range("A2:C9").Select()
$excel.Selection.EntireRow.Delete()

# --- Crops: delete data rows 2:5 (selection was A5, now A2:A5) ---
$ws = $wb.Worksheets.Item("Crops")
$ws.Activate()
$ws.Range("A2:A5").Select()
$excel.Selection.EntireRow.Delete()

# --- Researchers: delete data rows 2:6 (selection was A6, now A2:A6) ---
$ws = $wb.Worksheets.Item("Researchers")
$ws.Activate()
$ws.Range("A2:A6").Select()
$excel.Selection.EntireRow.Delete()

# --- MetStations: delete data rows 2:5 (selection was A6, now A2:A5) ---
$ws = $wb.Worksheets.Item("MetStations")
$ws.Activate()
$ws.Range("A2:A5").Select()
$excel.Selection.EntireRow.Delete()

# --- Fertilizers: delete data rows 2:5 (selection was H5, now A2:A5) ---
$ws = $wb.Worksheets.Item("Fertilizers")
$ws.Activate()
$ws.Range("A2:A5").Select()
$excel.Selection.EntireRow.Delete()

# --- Methods: delete data rows 2:5 (selection was A6, now A2:A5) ---
$ws = $wb.Worksheets.Item("Methods")
$ws.Activate()
$ws.Range("A2:A5").Select()
$excel.Selection.EntireRow.Delete()

# --- Units: delete data rows 2:6 (selection was A6, now A2:A6); ends up the active sheet ---
$ws = $wb.Worksheets.Item("Units")
$ws.Activate()
$ws.Range("A2:A6").Select()
$excel.Selection.EntireRow.Delete()
